# DBSettings.xlsx schema change: rebuild the "Data" sheet rows 2-10 with the
# new set of configuration rows, dropping the old NIKI*-related rows and the
# trailing rows 11-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete trailing rows (old rows 11-14) entirely so the sheet
# dimension shrinks back down to A1:H10.
$ws.Range("A11:H14").Delete()

# Column A (Id) holds numeric-looking text ("2", "3", ...). Force the whole
# range to Text format up-front so assigning the digit strings keeps them as
# shared-string text cells instead of being auto-converted to numbers.
$ws.Range("A2:A10").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C10").NumberFormat = "@"

# Row 2: OrgName
$ws.Range("A2").Value = "2"
$ws.Range("B2").Value = "OrgName"
$ws.Range("C2").Value = "КЛИНИКА"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Полное название ЛПУ"
$ws.Range("F2").Value = ""

# Row 3: NotificationServiceAddress
$ws.Range("A3").Value = "3"
$ws.Range("B3").Value = "NotificationServiceAddress"
$ws.Range("C3").Value = "net.tcp://localhost:8733/NotificationServiceEngine"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Адрес сервиса оповещений"
$ws.Range("F3").Value = ""

# Row 4: OrgOKPO
$ws.Range("A4").Value = "4"
$ws.Range("B4").Value = "OrgOKPO"
$ws.Range("C4").Value = "11223444"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "ОКПО ЛПУ"
$ws.Range("F4").Value = ""

# Row 5: OrgShortName
$ws.Range("A5").Value = "5"
$ws.Range("B5").Value = "OrgShortName"
$ws.Range("C5").Value = "КЛИНИКА"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Сокращенное название ЛПУ"
$ws.Range("F5").Value = ""

# Row 6: DirectorFullName
$ws.Range("A6").Value = "7"
$ws.Range("B6").Value = "DirectorFullName"
$ws.Range("C6").Value = "РУКОВОДИТЕЛЬ"
$ws.Range("D6").Value = "РУКОВОДИТЕЛЬ"
$ws.Range("E6").Value = "Руководитель"
$ws.Range("F6").Value = ""

# Row 7: PayContractLicense
$ws.Range("A7").Value = "8"
$ws.Range("B7").Value = "PayContractLicense"
$ws.Range("C7").Value = "ДОВЕРЕННОСТЬ"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "Доверенность на оказание платных услуг"
$ws.Range("F7").Value = ""

# Row 8: OrgAddress
$ws.Range("A8").Value = "10"
$ws.Range("B8").Value = "OrgAddress"
$ws.Range("C8").Value = "АДРЕС"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "Юридический адрес"
$ws.Range("F8").Value = ""

# Row 9: DirectorShortName
$ws.Range("A9").Value = "11"
$ws.Range("B9").Value = "DirectorShortName"
$ws.Range("C9").Value = "РУКОВОДИТЕЛЬ"
$ws.Range("D9").Value = "РУКОВОДИТЕЛЬ"
$ws.Range("E9").Value = "Руководитель"
$ws.Range("F9").Value = ""

# Row 10: ChildAge
$ws.Range("A10").Value = "13"
$ws.Range("B10").Value = "ChildAge"
$ws.Range("C10").Value = "15"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "Возраст, старше которого человек считается взрослым"
$ws.Range("F10").Value = ""
